$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4: No=3, Mood=5 (text), Description="excellent service", Verification="yes"
$ws.Range("A4").Value = 3
$ws.Range("B4:B5").NumberFormat = "@"
$ws.Range("B4").Value = "5"
$ws.Range("C4").Value = "excellent service"
$ws.Range("D4").Value = "yes"

# Row 5: No=4, Mood=5 (text), Description="nice service", Verification="yes"
$ws.Range("A5").Value = 4
$ws.Range("B5").Value = "5"
$ws.Range("C5").Value = "nice service"
$ws.Range("D5").Value = "yes"
